# Applies:
#  1) The table on slide 16 switches to table style {4360092B-B7E1-4817-A4C6-FC57F66DCD48}.
#  2) The presentation's theme colour scheme (ppt/theme/theme2.xml, the theme
#     actually wired to the one real Slide Master / Design in this deck) is
#     repainted from the "Integral" palette to the stock Office palette - i.e.
#     the Design gallery swap that the author made ends up producing the same
#     "Office Theme" colours on the live master that used to live on the
#     (COM-unreachable) Notes Master's theme part.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}
$tableShape.Table.ApplyStyle("{4360092B-B7E1-4817-A4C6-FC57F66DCD48}")

# --- 2. Theme colours ------------------------------------------------------
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office Theme") 12-slot colour scheme, in
# dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink order.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

$themeColors = $slide16.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $entry = $themeColors.Colors($i)
    $entry.RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}
